$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Query text used across the refreshed "startup" sheet. Single-quoted
# here-strings so backticks / brackets / quotes inside the Cypher text are
# kept 100% literal (no PowerShell expansion).
# ---------------------------------------------------------------------------

$casesQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
 WHERE ss.disease_subtype IN ["Adenocarcinoma"] and sf.grouped_recurrence_score IN ["0-5"] and d.tumor_size_group In ["Not Reported"]
return ss.study_subject_id as `Case ID`,
       p.program_acronym as `Program Code`,
        p.program_id as Program_ID,
       s.study_acronym as `Arm`,
       ss.disease_subtype as `Diagnosis`,
       sf.grouped_recurrence_score AS `Recurrence Score`,
       d.tumor_size_group AS `tumor_size`,
       d.er_status AS `ER Status`,
       d.pr_status AS `PR Status`,
       demo.age_at_index AS `Age (years)`,
demo.survival_time AS `Survival (days)`
'@

$statQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE ss.disease_subtype IN ["Adenocarcinoma"] and sf.grouped_recurrence_score IN ["0-5"] and d.tumor_size_group In ["Not Reported"]
WITH ss
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files
'@

$samplesQuery = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE ss.disease_subtype IN ["Adenocarcinoma"] and sf.grouped_recurrence_score IN ["0-5"]  and d.tumor_size_group In ["Not Reported"]
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`,
            lp.test_name as Platform
'@

$filesQuery = @'
MATCH (f:file)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
 WHERE ss.disease_subtype IN ["Adenocarcinoma"] and sf.grouped_recurrence_score IN ["0-5"] and d.tumor_size_group In ["Not Reported"]
RETURN  f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
    f.file_size AS `Size`,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`,
    ss.disease_subtype as `Diagnosis`
'@

$neo4jFile = "TC01_Bento_Filter_Diagnosis-Adenocarcinoma_Neo4jData.xlsx"
$webFile = "TC01_Bento_Filter_Diagnosis-Adenocarcinoma_WebData.xlsx"

# ---------------------------------------------------------------------------
# New tab labels and new query bodies first (this is the order new distinct
# strings are introduced into the workbook's shared-string table).
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("A4").Value = "FilesTab"
$ws.Range("B3").Value = $samplesQuery
$ws.Range("B2").Value = $casesQuery
$ws.Range("B4").Value = $filesQuery

# ---------------------------------------------------------------------------
# Re-stamp the (unchanged) StatQuery text into column C on all three rows.
# C2 already holds this exact text, so touch it via a transient value first
# -- otherwise a same-value write is a no-op and C2 would keep referencing
# its original shared-string slot instead of being re-anchored alongside
# the freshly written C3/C4 cells.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "placeholder"
$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery

# ---------------------------------------------------------------------------
# Remaining cells in each row (these all reuse strings already present in
# the workbook, so they don't add new shared-string entries).
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = $neo4jFile
$ws.Range("E2").Value = $webFile

$ws.Range("D3").Value = $neo4jFile
$ws.Range("E3").Value = $webFile

$ws.Range("D4").Value = $neo4jFile
$ws.Range("E4").Value = $webFile

# ---------------------------------------------------------------------------
# Formatting: wrap text on the query columns (B/C) for the two new rows,
# matching the existing row 2 style, and set the row heights to fit the
# new/updated content.
# ---------------------------------------------------------------------------
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true
$ws.Range("B4").WrapText = $true
$ws.Range("C4").WrapText = $true

$ws.Rows.Item(2).RowHeight = 333.5
$ws.Rows.Item(3).RowHeight = 391.5
$ws.Rows.Item(4).RowHeight = 290

# ---------------------------------------------------------------------------
# View state: selection moved down now that there are more rows.
# ---------------------------------------------------------------------------
$ws.Range("C9").Select()
